$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header fields for naming consistency
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update selection to match the post-edit state (K1:L1, active cell K1)
$ws.Range("K1:L1").Select()
